$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: opportunity 1329970 (Logistics Coordinator Intern @ HILTI Panama) ---
$ws.Range("A2").Value = "1329970"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1329970"
$ws.Range("C2").Value = "Logistics Coordinator Intern"
$ws.Range("F2").Value = "3 applicants"
$ws.Range("H2").Value = "HILTI Panama"

# --- Row 3: opportunity 1329856 (Infosys InStep - Global Internship Program) ---
$ws.Range("A3").Value = "1329856"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1329856"
$ws.Range("C3").Value = "Infosys InStep - Global Internship Program"
$ws.Range("D3").Value = "Bangalore, Karnataka, India"
$ws.Range("E3").Interior.Color = 65535
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "4 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Infosys Limited"

# --- Row 4: opportunity 1329697 (Interior architect @ FacePro) ---
$ws.Range("A4").Value = "1329697"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1329697"
$ws.Range("C4").Value = "Interior architect"
$ws.Range("D4").Value = "Sfax, Tunisie"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "FacePro"

# --- Row 5: opportunity 1325604 (International Business & Innovation Analyst) ---
$ws.Range("A5").Value = "1325604"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1325604"
$ws.Range("C5").Value = "International Business & Innovation Analyst"
$ws.Range("D5").Value = "4520 Santa Maria da Feira, Portugal"
$ws.Range("F5").Value = "165 applicants"
$ws.Range("H5").Value = "M2K Consultoria"

# --- Row 6: opportunity 1316099 (Sales & BD Junior Manager) ---
$ws.Range("A6").Value = "1316099"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1316099"
$ws.Range("C6").Value = "Sales & BD Junior Manager"
$ws.Range("D6").Value = "Puzi City, Chiayi County, Taiwan 613"
$ws.Range("F6").Value = "94 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Singform Enterprise Co., Ltd."

# --- Remove the old rows 7-12 (list shrank from 11 to 5 opportunities) ---
$ws.Rows("7:12").Delete()

# --- Narrower TITLE and ORGANIZATION columns ---
$ws.Columns("C").ColumnWidth = 45.083333333333336
$ws.Columns("H").ColumnWidth = 31.083333333333332
